# Apply the "Country" -> "Country qualifications" column header rename and
# update border-location rows that now qualify for multiple countries.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Linking station names")
$ws2 = $wb.Worksheets.Item("Station addresses")

# --- "Station addresses" sheet ---------------------------------------

# Border stations that should list both countries instead of a single one
$multiCountryRows = @(11, 12, 30, 41, 60, 78, 85)
foreach ($r in $multiCountryRows) {
    $ws2.Range("C$r").Value2 = "Netherlands, Germany"
}

# Column C header: "Country" -> "Country qualifications"
$ws2.Range("C1").Value2 = "Country qualifications"

# Widen column C to fit the new, longer values (target character width ~29.18)
$ws2.Columns("C").ColumnWidth = 28.4

# --- Cosmetic view-state updates (selection / active cell) ------------

$ws1.Activate()
$ws1.Range("C90").Select() | Out-Null

$ws2.Activate()
$ws2.Range("H13").Select() | Out-Null
